$d = $word.ActiveDocument

$replacements = @(
    @("159÷6=", "832÷3="),
    @("498÷7=", "469÷2="),
    @("150÷9=", "940÷8="),
    @("823÷7=", "977÷2="),
    @("923÷5=", "178÷8="),
    @("260÷2=", "100÷6="),
    @("246÷7=", "827÷5="),
    @("651÷7=", "513÷9="),
    @("724÷4=", "646÷7="),
    @("256÷4=", "923÷3="),
    @("116÷6=", "186÷5="),
    @("980÷3=", "990÷6="),
    @("892÷5=", "884÷4="),
    @("668÷4=", "530÷5="),
    @("746÷7=", "986÷9="),
    @("843÷4=", "926÷6="),
    @("825÷7=", "526÷2="),
    @("360÷4=", "523÷2="),
    @("661÷4=", "844÷8="),
    @("526÷5=", "223÷2="),
    @("791÷8=", "133÷5="),
    @("278÷2=", "382÷8="),
    @("262÷6=", "135÷4="),
    @("439÷7=", "608÷5="),
    @("107÷6=", "995÷3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
